$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 19787.5
$ws.Range("I18").Value = 25980
$ws.Range("J18").Value = 9466.666999999999
$ws.Range("K18").Value = 25980
$ws.Range("L18").Value = 9466.666999999999
$ws.Range("M18").Value = -25696
$ws.Range("N18").Value = -10034.667
$ws.Range("H43").Value = 3598.3333
$ws.Range("I43").Value = 3995
$ws.Range("J43").Value = 3400
$ws.Range("K43").Value = 3995
$ws.Range("L43").Value = 3400
$ws.Range("M43").Value = -3926
$ws.Range("N43").Value = -3538
$ws.Range("H86").Value = 2000
$ws.Range("J86").Value = 2000
$ws.Range("L86").Value = 2000
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 2000
$ws.Range("J89").Value = 2000
$ws.Range("L89").Value = 10000
$ws.Range("N89").Value = -21232
$ws.Range("H132").Value = 2627.5667
$ws.Range("I132").Value = 2373.5518
$ws.Range("K132").Value = 7120.655400000001
$ws.Range("M132").Value = -4590.655400000001
$ws.Range("H137").Value = 2408.3333
$ws.Range("I137").Value = 2408.3333
$ws.Range("K137").Value = 7224.999899999999
$ws.Range("M137").Value = -4674.999899999999
$ws.Range("H138").Value = 3297.5
$ws.Range("J138").Value = 3925.7144
$ws.Range("L138").Value = 11777.1432
$ws.Range("N138").Value = -22057.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H106").Value = 40000
$ws.Range("J106").Value = 40000
$ws.Range("L106").Value = 40000
$ws.Range("N106").Value = -42524

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 488.33334
$ws.Range("I80").Value = 319.2
$ws.Range("J80").Value = 699.75
$ws.Range("K80").Value = 319.2
$ws.Range("L80").Value = 699.75
$ws.Range("M80").Value = 678.8
$ws.Range("N80").Value = -2695.75
$ws.Range("H83").Value = 488.33334
$ws.Range("I83").Value = 319.2
$ws.Range("J83").Value = 699.75
$ws.Range("K83").Value = 1596
$ws.Range("L83").Value = 3498.75
$ws.Range("M83").Value = 3396
$ws.Range("N83").Value = -13482.75
$ws.Range("H86").Value = 25834.5
$ws.Range("I86").Value = 1500
$ws.Range("J86").Value = 38001.75
$ws.Range("K86").Value = 1500
$ws.Range("L86").Value = 38001.75
$ws.Range("M86").Value = -377
$ws.Range("N86").Value = -40247.75
$ws.Range("H89").Value = 25834.5
$ws.Range("I89").Value = 1500
$ws.Range("J89").Value = 38001.75
$ws.Range("K89").Value = 7500
$ws.Range("L89").Value = 190008.75
$ws.Range("M89").Value = -1884
$ws.Range("N89").Value = -201240.75
$ws.Range("H94").Value = 3768.5
$ws.Range("I94").Value = 2531.2
$ws.Range("K94").Value = 2531.2
$ws.Range("M94").Value = -2080.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H13").Value = 1930.6
$ws.Range("I13").Value = 1634.3334
$ws.Range("J13").Value = 2375
$ws.Range("K13").Value = 1634.3334
$ws.Range("L13").Value = 2375
$ws.Range("M13").Value = -1495.3334
$ws.Range("N13").Value = -2653
$ws.Range("H20").Value = 99497
$ws.Range("J20").Value = 99497
$ws.Range("L20").Value = 99497
$ws.Range("N20").Value = -99969
$ws.Range("H30").Value = 99497
$ws.Range("J30").Value = 99497
$ws.Range("L30").Value = 99497
$ws.Range("N30").Value = -99679
$ws.Range("H31").Value = 3049.8
$ws.Range("J31").Value = 4952.778
$ws.Range("L31").Value = 4952.778
$ws.Range("N31").Value = -5542.778
$ws.Range("H34").Value = 3049.8
$ws.Range("J34").Value = 4952.778
$ws.Range("L34").Value = 4952.778
$ws.Range("N34").Value = -5356.778
$ws.Range("H62").Value = 15002.5
$ws.Range("I62").Value = 15002.5
$ws.Range("K62").Value = 15002.5
$ws.Range("M62").Value = -14378.5
$ws.Range("H65").Value = 15002.5
$ws.Range("I65").Value = 15002.5
$ws.Range("K65").Value = 75012.5
$ws.Range("M65").Value = -71892.5
$ws.Range("H128").Value = 99497
$ws.Range("J128").Value = 99497
$ws.Range("L128").Value = 99497
$ws.Range("N128").Value = -109457
$ws.Range("H134").Value = 2015.5834
$ws.Range("I134").Value = 2118.7
$ws.Range("K134").Value = 6356.099999999999
$ws.Range("M134").Value = -3821.099999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()
$ws.Range("H55").Value = 3074.0908
$ws.Range("J55").Value = 3335
$ws.Range("L55").Value = 10005
$ws.Range("N55").Value = -10359
$ws.Range("H131").Value = 799
$ws.Range("I131").Value = 799
$ws.Range("K131").Value = 2397
$ws.Range("M131").Value = 2643
$ws.Range("H136").Value = 994.4286
$ws.Range("I136").Value = 994.4286
$ws.Range("K136").Value = 2983.2858
$ws.Range("M136").Value = 2116.7142
$ws.Range("H138").Value = 2464.1428
$ws.Range("I138").Value = 2374.8333
$ws.Range("K138").Value = 7124.499899999999
$ws.Range("M138").Value = -1984.499899999999
$ws.Range("H139").Value = 2257.25
$ws.Range("I139").Value = 2257.25
$ws.Range("K139").Value = 6771.75
$ws.Range("M139").Value = -1631.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 101503
$ws.Range("I80").Value = 3000
$ws.Range("K80").Value = 3000
$ws.Range("M80").Value = -2002
$ws.Range("H83").Value = 101503
$ws.Range("I83").Value = 3000
$ws.Range("K83").Value = 15000
$ws.Range("M83").Value = -10008
$ws.Range("H95").Value = 49111.75
$ws.Range("J95").Value = 49111.75
$ws.Range("L95").Value = 49111.75
$ws.Range("N95").Value = -54603.75
$ws.Range("H132").Value = 3181.25
$ws.Range("I132").Value = 2658.3333
$ws.Range("K132").Value = 7974.999899999999
$ws.Range("M132").Value = -5444.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H25").Value = 3338.3333
$ws.Range("I25").Value = 2507.5
$ws.Range("J25").Value = 5000
$ws.Range("K25").Value = 2507.5
$ws.Range("L25").Value = 5000
$ws.Range("M25").Value = -2277.5
$ws.Range("N25").Value = -5460
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H36").Value = 99994.5
$ws.Range("J36").Value = 99994.5
$ws.Range("L36").Value = 99994.5
$ws.Range("N36").Value = -101118.5
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()
$ws.Range("H46").Value = 3055.389
$ws.Range("I46").Value = 1500
$ws.Range("J46").Value = 3366.4666
$ws.Range("K46").Value = 1500
$ws.Range("L46").Value = 3366.4666
$ws.Range("M46").Value = -1312
$ws.Range("N46").Value = -3742.4666
$ws.Range("H68").Value = 2619.9
$ws.Range("I68").Value = 2577.6667
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 2577.6667
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -1828.6667
$ws.Range("N68").Value = -4498
$ws.Range("H71").Value = 2619.9
$ws.Range("I71").Value = 2577.6667
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 12888.3335
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -9144.333500000001
$ws.Range("N71").Value = -22488
$ws.Range("H82").Value = 2491.6667
$ws.Range("I82").Value = 2237.5
$ws.Range("J82").Value = 3000
$ws.Range("K82").Value = 2237.5
$ws.Range("L82").Value = 3000
$ws.Range("M82").Value = -1876.5
$ws.Range("N82").Value = -3722
$ws.Range("H85").Value = 2491.6667
$ws.Range("I85").Value = 2237.5
$ws.Range("J85").Value = 3000
$ws.Range("K85").Value = 2237.5
$ws.Range("L85").Value = 3000
$ws.Range("M85").Value = -989.5
$ws.Range("N85").Value = -5496
$ws.Range("H93").Value = 1914.2858
$ws.Range("J93").Value = 1925.5
$ws.Range("L93").Value = 1925.5
$ws.Range("N93").Value = -4421.5
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()
$ws.Range("H136").Value = 3640.8
$ws.Range("I136").Value = 3640.8
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10922.4
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -8372.400000000001
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 3709.3333
$ws.Range("I3").Value = 564.5
$ws.Range("J3").Value = 9999
$ws.Range("K3").Value = 564.5
$ws.Range("L3").Value = 9999
$ws.Range("M3").Value = -450.5
$ws.Range("N3").Value = -10227
$ws.Range("H21").Value = 3765007.5
$ws.Range("I21").Value = 3765007.5
$ws.Range("K21").Value = 3765007.5
$ws.Range("M21").Value = -3764772.5
$ws.Range("H35").Value = 3765007.5
$ws.Range("I35").Value = 3765007.5
$ws.Range("K35").Value = 3765007.5
$ws.Range("M35").Value = -3764717.5
$ws.Range("H58").Value = 28047.5
$ws.Range("I58").Value = 30028.334
$ws.Range("K58").Value = 30028.334
$ws.Range("M58").Value = -29720.334
$ws.Range("H74").Value = 15000
$ws.Range("J74").Value = 15000
$ws.Range("L74").Value = 15000
$ws.Range("N74").Value = -16872
$ws.Range("H77").Value = 15000
$ws.Range("J77").Value = 15000
$ws.Range("L77").Value = 45000
$ws.Range("N77").Value = -54360
$ws.Range("H135").Value = 99400
$ws.Range("J135").Value = 99400
$ws.Range("L135").Value = 99400
$ws.Range("N135").Value = -109540
